$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6
$ws.Range("A6").Value = "REV-GV-300"
$ws.Range("B6").Value = "Customisable slope range"
$ws.Range("C6").Value = "implement customisable slope range in pattern widget"

# New row 7, column A
$ws.Range("A7").Value = "DEV-GV-300"

# Row 3, column A: Ref changes from DEV-GV-100 to DEV-GV-200
$ws.Range("A3").Value = "DEV-GV-200"

# New row 7, columns B and C
$ws.Range("B7").Value = "Extend pattern"
$ws.Range("C7").Value = "opposite function to shrink pattern"

# Update selection to match target state
$ws.Range("C7").Select()
